$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = 2
    4  = -1
    5  = -1
    6  = -1
    7  = -4
    8  = -3
    10 = 3
    11 = -1
    12 = -4
    13 = 0
    14 = 1
    15 = 2
    16 = -3
    17 = -4
    18 = 3
    19 = 1
    20 = -2
    22 = -1
    23 = -5
    24 = -1
    25 = -3
    26 = 4
    27 = 5
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
